# [EI-979] [Rollback] Changed "then_goto" and "else_goto" back to
# "then_question" and "else_question" in the survey.xlsx data dictionary.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the two header cells (I1/J1) that hold the column titles.
$ws.Range("I1").Value = "Then_Goto"
$ws.Range("J1").Value = "Else_Goto"

# Restore the sheet's active selection to I1 (matches the saved view state).
$ws.Range("I1").Select()
